# Applies the edits described in the commit:
#   "add animation for chatbot, make new error log for windows error"
#
# The sheet is a flat operation/error log table with columns:
#   A=time B=type C=user_name D=pc_name E=win_title F=win_urlpath
#   G=win_hwnd H=win_class I=app_path J=capimg K=explanation
#   L=error_type M=error_content
#
# The edit:
#  - renames the user on every data row (2-16) from "Haruka Sasaki" to
#    "Akira Kobayashi"
#  - renumbers/renames the screenshot file names in column J
#  - rewrites most of the step explanations in column K (the steps were
#    reshuffled / reworded to describe a slightly different operation
#    sequence, and a new Windows-update error, 0x80240fff, replaces the
#    old 0x80244007 one, shifting which row is flagged as an "error")
#  - row 5 becomes the "error" row (0x80240fff) and row 7 becomes a
#    normal "operation" row again (the old error row 7, 0x80244007, is
#    cleared out)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- column B (type) -------------------------------------------------
$ws.Range("B5").Value = "error"
$ws.Range("B7").Value = "operation"

# ---- column C (user_name) : same rename applied to rows 2-16 --------
foreach ($r in 2..16) {
    $ws.Cells.Item($r, 3).Value = "Akira Kobayashi"
}

# ---- column J (capimg) ------------------------------------------------
$ws.Range("J2").Value  = "bdot20240415_141954/1.png"
$ws.Range("J3").Value  = "bdot20240415_141954/2.png"
$ws.Range("J4").Value  = "bdot20240415_141954/3.png"
$ws.Range("J5").Value  = "bdot20240415_141954/4.png"
$ws.Range("J6").Value  = "bdot20240415_141954/5.png"
$ws.Range("J7").Value  = "bdot20240415_141954/5.png"
$ws.Range("J8").Value  = "bdot20240415_141954/6.png"
$ws.Range("J9").Value  = "bdot20240415_141954/7.png"
$ws.Range("J10").Value = "bdot20240415_141954/8.png"
$ws.Range("J11").Value = "bdot20240415_141954/9.png"
$ws.Range("J12").Value = "bdot20240415_141954/10.png"
$ws.Range("J13").Value = "bdot20240415_141954/1.png"
$ws.Range("J14").Value = "bdot20240415_141954/2.png"
$ws.Range("J15").Value = "bdot20240415_141954/3.png"
$ws.Range("J16").Value = "bdot20240415_141954/11.png"

# ---- column K (explanation) -------------------------------------------
$ws.Range("K2").Value  = "「スタート」ボタンをクリックする"
$ws.Range("K3").Value  = "メニューから「設定」アイコンをクリックする"
$ws.Range("K4").Value  = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K5").Value  = "0x80240fff エラー"
$ws.Range("K6").Value  = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
$ws.Range("K7").Value  = "メニューからターミナル(管理者)をクリックする"
$ws.Range("K8").Value  = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
$ws.Range("K9").Value  = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
$ws.Range("K10").Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
$ws.Range("K11").Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"
$ws.Range("K12").Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
$ws.Range("K13").Value = "「スタート」ボタンをクリックする"
$ws.Range("K14").Value = "メニューから「設定」アイコンをクリックする"
$ws.Range("K15").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K16").Value = "「更新プログラムのチェック」ボタンをクリックする"

# ---- columns L/M (error_type / error_content) -------------------------
# Row 5 now carries the new Windows-update error details.
$ws.Range("L5").Value = "Error W"
$ws.Range("M5").Value = " エラーの Windows"

# Row 7 is no longer an error row, so its error columns are cleared.
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
